$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 234,5
$arr[0,0] = [double]"6"
$arr[0,1] = "CT"
$arr[0,2] = [double]"1"
$arr[0,3] = [double]"413071"
$arr[0,4] = [double]"1.1609237009213569e-11"
$arr[1,0] = [double]"7"
$arr[1,1] = "CT"
$arr[1,2] = [double]"1"
$arr[1,3] = [double]"413071"
$arr[1,4] = [double]"9.4581330614684767e-13"
$arr[2,0] = [double]"8"
$arr[2,1] = "CT"
$arr[2,2] = [double]"7382.95458984375"
$arr[2,3] = [double]"413071"
$arr[2,4] = [double]"5.8013052139926913e-09"
$arr[3,0] = [double]"9"
$arr[3,1] = "CT"
$arr[3,2] = [double]"9492.3701171875"
$arr[3,3] = [double]"413071"
$arr[3,4] = [double]"6.5862977294273151e-09"
$arr[4,0] = [double]"10"
$arr[4,1] = "CT"
$arr[4,2] = [double]"21094.158203125"
$arr[4,3] = [double]"413071"
$arr[4,4] = [double]"3.7394848106941936e-08"
$arr[5,0] = [double]"11"
$arr[5,1] = "CT"
$arr[5,2] = [double]"15820.6162109375"
$arr[5,3] = [double]"413071"
$arr[5,4] = [double]"2.9181531857602749e-08"
$arr[6,0] = [double]"12"
$arr[6,1] = "CT"
$arr[6,2] = [double]"52735.390625"
$arr[6,3] = [double]"413071"
$arr[6,4] = [double]"1.3141298893515341e-07"
$arr[7,0] = [double]"13"
$arr[7,1] = "CT"
$arr[7,2] = [double]"132927.171875"
$arr[7,3] = [double]"413071"
$arr[7,4] = [double]"2.4623619765407057e-07"
$arr[8,0] = [double]"14"
$arr[8,1] = "CT"
$arr[8,2] = [double]"137798.328125"
$arr[8,3] = [double]"413071"
$arr[8,4] = [double]"1.9000468398644443e-07"
$arr[9,0] = [double]"15"
$arr[9,1] = "CT"
$arr[9,2] = [double]"136921.75"
$arr[9,3] = [double]"413071"
$arr[9,4] = [double]"1.0870907374282979e-07"
$arr[10,0] = [double]"16"
$arr[10,1] = "CT"
$arr[10,2] = [double]"190174.25"
$arr[10,3] = [double]"413071"
$arr[10,4] = [double]"3.0874960543769703e-07"
$arr[11,0] = [double]"17"
$arr[11,1] = "CT"
$arr[11,2] = [double]"197436.4375"
$arr[11,3] = [double]"413071"
$arr[11,4] = [double]"3.5174818435734778e-07"
$arr[12,0] = [double]"18"
$arr[12,1] = "CT"
$arr[12,2] = [double]"161044.625"
$arr[12,3] = [double]"413071"
$arr[12,4] = [double]"5.0302753606956685e-07"
$arr[13,0] = [double]"19"
$arr[13,1] = "CT"
$arr[13,2] = [double]"41829.15625"
$arr[13,3] = [double]"413071"
$arr[13,4] = [double]"1.3645779972648597e-07"
$arr[14,0] = [double]"20"
$arr[14,1] = "CT"
$arr[14,2] = [double]"17886.884765625"
$arr[14,3] = [double]"413071"
$arr[14,4] = [double]"7.1206400775736256e-08"
$arr[15,0] = [double]"21"
$arr[15,1] = "CT"
$arr[15,2] = [double]"14599.7353515625"
$arr[15,3] = [double]"413071"
$arr[15,4] = [double]"1.0856475540776955e-07"
$arr[16,0] = [double]"22"
$arr[16,1] = "CT"
$arr[16,2] = [double]"8551.65234375"
$arr[16,3] = [double]"413071"
$arr[16,4] = [double]"4.5457515085445266e-08"
$arr[17,0] = [double]"23"
$arr[17,1] = "CT"
$arr[17,2] = [double]"4908.30029296875"
$arr[17,3] = [double]"413071"
$arr[17,4] = [double]"3.6734391528625565e-08"
$arr[18,0] = [double]"24"
$arr[18,1] = "CT"
$arr[18,2] = [double]"3174.176025390625"
$arr[18,3] = [double]"413071"
$arr[18,4] = [double]"2.3057520337488313e-08"
$arr[19,0] = [double]"25"
$arr[19,1] = "CT"
$arr[19,2] = [double]"542.42962646484375"
$arr[19,3] = [double]"413071"
$arr[19,4] = [double]"8.635289816538716e-09"
$arr[20,0] = [double]"26"
$arr[20,1] = "CT"
$arr[20,2] = [double]"129.90037536621094"
$arr[20,3] = [double]"413071"
$arr[20,4] = [double]"2.0487129770430101e-09"
$arr[21,0] = [double]"27"
$arr[21,1] = "CT"
$arr[21,2] = [double]"293.73092651367188"
$arr[21,3] = [double]"413071"
$arr[21,4] = [double]"1.008767025467705e-08"
$arr[22,0] = [double]"28"
$arr[22,1] = "CT"
$arr[22,2] = [double]"469.26144409179688"
$arr[22,3] = [double]"413071"
$arr[22,4] = [double]"1.5916977957886047e-08"
$arr[23,0] = [double]"29"
$arr[23,1] = "CT"
$arr[23,2] = [double]"1"
$arr[23,3] = [double]"413071"
$arr[23,4] = [double]"9.9037306133809011e-11"
$arr[24,0] = [double]"30"
$arr[24,1] = "CT"
$arr[24,2] = [double]"34.604843139648438"
$arr[24,3] = [double]"413071"
$arr[24,4] = [double]"8.7307503449096657e-09"
$arr[25,0] = [double]"31"
$arr[25,1] = "CT"
$arr[25,2] = [double]"1"
$arr[25,3] = [double]"413071"
$arr[25,4] = [double]"3.7844777311946132e-10"
$arr[26,0] = [double]"6"
$arr[26,1] = "DE"
$arr[26,2] = [double]"1"
$arr[26,3] = [double]"401288"
$arr[26,4] = [double]"5.5119631081623766e-11"
$arr[27,0] = [double]"7"
$arr[27,1] = "DE"
$arr[27,2] = [double]"3940.27587890625"
$arr[27,3] = [double]"401288"
$arr[27,4] = [double]"1.7694352649755274e-08"
$arr[28,0] = [double]"8"
$arr[28,1] = "DE"
$arr[28,2] = [double]"3979.85302734375"
$arr[28,3] = [double]"401288"
$arr[28,4] = [double]"1.4847905127624017e-08"
$arr[29,0] = [double]"9"
$arr[29,1] = "DE"
$arr[29,2] = [double]"5969.779296875"
$arr[29,3] = [double]"401288"
$arr[29,4] = [double]"1.9666524408989972e-08"
$arr[30,0] = [double]"10"
$arr[30,1] = "DE"
$arr[30,2] = [double]"11939.8681640625"
$arr[30,3] = [double]"401288"
$arr[30,4] = [double]"1.0049668475176077e-07"
$arr[31,0] = [double]"11"
$arr[31,1] = "DE"
$arr[31,2] = [double]"23879.921875"
$arr[31,3] = [double]"401288"
$arr[31,4] = [double]"2.0913186915549886e-07"
$arr[32,0] = [double]"12"
$arr[32,1] = "DE"
$arr[32,2] = [double]"45768.3125"
$arr[32,3] = [double]"401288"
$arr[32,4] = [double]"5.4150649475559476e-07"
$arr[33,0] = [double]"13"
$arr[33,1] = "DE"
$arr[33,2] = [double]"77610.2578125"
$arr[33,3] = [double]"401288"
$arr[33,4] = [double]"6.8258992769187898e-07"
$arr[34,0] = [double]"14"
$arr[34,1] = "DE"
$arr[34,2] = [double]"117786.9765625"
$arr[34,3] = [double]"401288"
$arr[34,4] = [double]"7.7111695873099961e-07"
$arr[35,0] = [double]"15"
$arr[35,1] = "DE"
$arr[35,2] = [double]"118327.9453125"
$arr[35,3] = [double]"401288"
$arr[35,4] = [double]"4.4604976778828132e-07"
$arr[36,0] = [double]"16"
$arr[36,1] = "DE"
$arr[36,2] = [double]"56474.09765625"
$arr[36,3] = [double]"401288"
$arr[36,4] = [double]"4.3531801452445507e-07"
$arr[37,0] = [double]"17"
$arr[37,1] = "DE"
$arr[37,2] = [double]"26512.328125"
$arr[37,3] = [double]"401288"
$arr[37,4] = [double]"2.242616687908594e-07"
$arr[38,0] = [double]"18"
$arr[38,1] = "DE"
$arr[38,2] = [double]"19543.951171875"
$arr[38,3] = [double]"401288"
$arr[38,4] = [double]"2.8984149480493215e-07"
$arr[39,0] = [double]"19"
$arr[39,1] = "DE"
$arr[39,2] = [double]"10862.0146484375"
$arr[39,3] = [double]"401288"
$arr[39,4] = [double]"1.6824115789404459e-07"
$arr[40,0] = [double]"20"
$arr[40,1] = "DE"
$arr[40,2] = [double]"3729.3203125"
$arr[40,3] = [double]"401288"
$arr[40,4] = [double]"7.048823391642145e-08"
$arr[41,0] = [double]"21"
$arr[41,1] = "DE"
$arr[41,2] = [double]"1345.5067138671875"
$arr[41,3] = [double]"401288"
$arr[41,4] = [double]"4.7504240541229592e-08"
$arr[42,0] = [double]"22"
$arr[42,1] = "DE"
$arr[42,2] = [double]"3192.2021484375"
$arr[42,3] = [double]"401288"
$arr[42,4] = [double]"8.056542100121078e-08"
$arr[43,0] = [double]"23"
$arr[43,1] = "DE"
$arr[43,2] = [double]"25.691064834594727"
$arr[43,3] = [double]"401288"
$arr[43,4] = [double]"9.1290680570565996e-10"
$arr[44,0] = [double]"24"
$arr[44,1] = "DE"
$arr[44,2] = [double]"254.52424621582031"
$arr[44,3] = [double]"401288"
$arr[44,4] = [double]"8.7783602609192712e-09"
$arr[45,0] = [double]"25"
$arr[45,1] = "DE"
$arr[45,2] = [double]"178.62748718261719"
$arr[45,3] = [double]"401288"
$arr[45,4] = [double]"1.3501574080976297e-08"
$arr[46,0] = [double]"26"
$arr[46,1] = "DE"
$arr[46,2] = [double]"881.1624755859375"
$arr[46,3] = [double]"401288"
$arr[46,4] = [double]"6.5982590058411006e-08"
$arr[47,0] = [double]"27"
$arr[47,1] = "DE"
$arr[47,2] = [double]"1"
$arr[47,3] = [double]"401288"
$arr[47,4] = [double]"1.6305864991572605e-10"
$arr[48,0] = [double]"28"
$arr[48,1] = "DE"
$arr[48,2] = [double]"1"
$arr[48,3] = [double]"401288"
$arr[48,4] = [double]"1.6104542699402202e-10"
$arr[49,0] = [double]"29"
$arr[49,1] = "DE"
$arr[49,2] = [double]"1"
$arr[49,3] = [double]"401288"
$arr[49,4] = [double]"4.7022036264721123e-10"
$arr[50,0] = [double]"30"
$arr[50,1] = "DE"
$arr[50,2] = [double]"1"
$arr[50,3] = [double]"401288"
$arr[50,4] = [double]"1.1978910086085648e-09"
$arr[51,0] = [double]"31"
$arr[51,1] = "DE"
$arr[51,2] = [double]"1"
$arr[51,3] = [double]"401288"
$arr[51,4] = [double]"1.7968365684239984e-09"
$arr[52,0] = [double]"6"
$arr[52,1] = "MA"
$arr[52,2] = [double]"1"
$arr[52,3] = [double]"104251"
$arr[52,4] = [double]"4.5998975461980507e-11"
$arr[53,0] = [double]"7"
$arr[53,1] = "MA"
$arr[53,2] = [double]"1"
$arr[53,3] = [double]"104251"
$arr[53,4] = [double]"3.7475712637291547e-12"
$arr[54,0] = [double]"8"
$arr[54,1] = "MA"
$arr[54,2] = [double]"1"
$arr[54,3] = [double]"104251"
$arr[54,4] = [double]"3.1134363750323635e-12"
$arr[55,0] = [double]"9"
$arr[55,1] = "MA"
$arr[55,2] = [double]"1"
$arr[55,3] = [double]"104251"
$arr[55,4] = [double]"2.7492307475701638e-12"
$arr[56,0] = [double]"10"
$arr[56,1] = "MA"
$arr[56,2] = [double]"9238.4013671875"
$arr[56,3] = [double]"104251"
$arr[56,4] = [double]"6.4891956697010755e-08"
$arr[57,0] = [double]"11"
$arr[57,1] = "MA"
$arr[57,2] = [double]"9238.4013671875"
$arr[57,3] = [double]"104251"
$arr[57,4] = [double]"6.7518996615945071e-08"
$arr[58,0] = [double]"12"
$arr[58,1] = "MA"
$arr[58,2] = [double]"18476.802734375"
$arr[58,3] = [double]"104251"
$arr[58,4] = [double]"1.8243468957734876e-07"
$arr[59,0] = [double]"13"
$arr[59,1] = "MA"
$arr[59,2] = [double]"37430.94921875"
$arr[59,3] = [double]"104251"
$arr[59,4] = [double]"2.7473464569993666e-07"
$arr[60,0] = [double]"14"
$arr[60,1] = "MA"
$arr[60,2] = [double]"27100.591796875"
$arr[60,3] = [double]"104251"
$arr[60,4] = [double]"1.4806199999384262e-07"
$arr[61,0] = [double]"15"
$arr[61,1] = "MA"
$arr[61,2] = [double]"28117.837890625"
$arr[61,3] = [double]"104251"
$arr[61,4] = [double]"8.845446330951745e-08"
$arr[62,0] = [double]"16"
$arr[62,1] = "MA"
$arr[62,2] = [double]"97040.53125"
$arr[62,3] = [double]"104251"
$arr[62,4] = [double]"6.2424106772596133e-07"
$arr[63,0] = [double]"17"
$arr[63,1] = "MA"
$arr[63,2] = [double]"16215.6904296875"
$arr[63,3] = [double]"104251"
$arr[63,4] = [double]"1.1446810077586633e-07"
$arr[64,0] = [double]"18"
$arr[64,1] = "MA"
$arr[64,2] = [double]"10365.4443359375"
$arr[64,3] = [double]"104251"
$arr[64,4] = [double]"1.282856061379789e-07"
$arr[65,0] = [double]"19"
$arr[65,1] = "MA"
$arr[65,2] = [double]"5789.8037109375"
$arr[65,3] = [double]"104251"
$arr[65,4] = [double]"7.4838929720044689e-08"
$arr[66,0] = [double]"20"
$arr[66,1] = "MA"
$arr[66,2] = [double]"10714.201171875"
$arr[66,3] = [double]"104251"
$arr[66,4] = [double]"1.6900074228942685e-07"
$arr[67,0] = [double]"21"
$arr[67,1] = "MA"
$arr[67,2] = [double]"3806.5068359375"
$arr[67,3] = [double]"104251"
$arr[67,4] = [double]"1.1215405493203434e-07"
$arr[68,0] = [double]"22"
$arr[68,1] = "MA"
$arr[68,2] = [double]"3239.448974609375"
$arr[68,3] = [double]"104251"
$arr[68,4] = [double]"6.8229354610593873e-08"
$arr[69,0] = [double]"23"
$arr[69,1] = "MA"
$arr[69,2] = [double]"2481.369873046875"
$arr[69,3] = [double]"104251"
$arr[69,4] = [double]"7.3583031223733997e-08"
$arr[70,0] = [double]"24"
$arr[70,1] = "MA"
$arr[70,2] = [double]"497.12643432617188"
$arr[70,3] = [double]"104251"
$arr[70,4] = [double]"1.4308461082634949e-08"
$arr[71,0] = [double]"25"
$arr[71,1] = "MA"
$arr[71,2] = [double]"1"
$arr[71,3] = [double]"104251"
$arr[71,4] = [double]"6.3078008094574756e-11"
$arr[72,0] = [double]"26"
$arr[72,1] = "MA"
$arr[72,2] = [double]"1"
$arr[72,3] = [double]"104251"
$arr[72,4] = [double]"6.2490672358972432e-11"
$arr[73,0] = [double]"27"
$arr[73,1] = "MA"
$arr[73,2] = [double]"53.702713012695313"
$arr[73,3] = [double]"104251"
$arr[73,4] = [double]"7.3077202067395319e-09"
$arr[74,0] = [double]"28"
$arr[74,1] = "MA"
$arr[74,2] = [double]"1"
$arr[74,3] = [double]"104251"
$arr[74,4] = [double]"1.3439720170094205e-10"
$arr[75,0] = [double]"29"
$arr[75,1] = "MA"
$arr[75,2] = [double]"1"
$arr[75,3] = [double]"104251"
$arr[75,4] = [double]"3.9241293547931377e-10"
$arr[76,0] = [double]"30"
$arr[76,1] = "MA"
$arr[76,2] = [double]"1"
$arr[76,3] = [double]"104251"
$arr[76,4] = [double]"9.9967567557257553e-10"
$arr[77,0] = [double]"31"
$arr[77,1] = "MA"
$arr[77,2] = [double]"1"
$arr[77,3] = [double]"104251"
$arr[77,4] = [double]"1.4995137354034682e-09"
$arr[78,0] = [double]"6"
$arr[78,1] = "MD"
$arr[78,2] = [double]"29103.9921875"
$arr[78,3] = [double]"544681"
$arr[78,4] = [double]"1.1818784741990385e-06"
$arr[79,0] = [double]"7"
$arr[79,1] = "MD"
$arr[79,2] = [double]"50032.640625"
$arr[79,3] = [double]"544681"
$arr[79,4] = [double]"1.6552951365156332e-07"
$arr[80,0] = [double]"8"
$arr[80,1] = "MD"
$arr[80,2] = [double]"55395.890625"
$arr[80,3] = [double]"544681"
$arr[80,4] = [double]"1.5226133598389424e-07"
$arr[81,0] = [double]"9"
$arr[81,1] = "MD"
$arr[81,2] = [double]"71883.3828125"
$arr[81,3] = [double]"544681"
$arr[81,4] = [double]"1.7446639333229541e-07"
$arr[82,0] = [double]"10"
$arr[82,1] = "MD"
$arr[82,2] = [double]"128631.1328125"
$arr[82,3] = [double]"544681"
$arr[82,4] = [double]"7.9764976135265897e-07"
$arr[83,0] = [double]"11"
$arr[83,1] = "MD"
$arr[83,2] = [double]"223461.390625"
$arr[83,3] = [double]"544681"
$arr[83,4] = [double]"1.4417958027479472e-06"
$arr[84,0] = [double]"12"
$arr[84,1] = "MD"
$arr[84,2] = [double]"140919.203125"
$arr[84,3] = [double]"544681"
$arr[84,4] = [double]"1.228352061843907e-06"
$arr[85,0] = [double]"13"
$arr[85,1] = "MD"
$arr[85,2] = [double]"76175.546875"
$arr[85,3] = [double]"544681"
$arr[85,4] = [double]"4.9359442755303462e-07"
$arr[86,0] = [double]"14"
$arr[86,1] = "MD"
$arr[86,2] = [double]"61683.2578125"
$arr[86,3] = [double]"544681"
$arr[86,4] = [double]"2.9751183205917187e-07"
$arr[87,0] = [double]"15"
$arr[87,1] = "MD"
$arr[87,2] = [double]"65187.88671875"
$arr[87,3] = [double]"544681"
$arr[87,4] = [double]"1.8104096000115533e-07"
$arr[88,0] = [double]"16"
$arr[88,1] = "MD"
$arr[88,2] = [double]"33684.73828125"
$arr[88,3] = [double]"544681"
$arr[88,4] = [double]"1.9129537065509794e-07"
$arr[89,0] = [double]"17"
$arr[89,1] = "MD"
$arr[89,2] = [double]"33559.91015625"
$arr[89,3] = [double]"544681"
$arr[89,4] = [double]"2.0914231413371454e-07"
$arr[90,0] = [double]"18"
$arr[90,1] = "MD"
$arr[90,2] = [double]"17296.1953125"
$arr[90,3] = [double]"544681"
$arr[90,4] = [double]"1.8897863185429742e-07"
$arr[91,0] = [double]"19"
$arr[91,1] = "MD"
$arr[91,2] = [double]"17406.576171875"
$arr[91,3] = [double]"544681"
$arr[91,4] = [double]"1.9863198019720585e-07"
$arr[92,0] = [double]"20"
$arr[92,1] = "MD"
$arr[92,2] = [double]"8982.533203125"
$arr[92,3] = [double]"544681"
$arr[92,4] = [double]"1.2508341740158357e-07"
$arr[93,0] = [double]"21"
$arr[93,1] = "MD"
$arr[93,2] = [double]"1941.474365234375"
$arr[93,3] = [double]"544681"
$arr[93,4] = [double]"5.0500084114446508e-08"
$arr[94,0] = [double]"22"
$arr[94,1] = "MD"
$arr[94,2] = [double]"1389.294921875"
$arr[94,3] = [double]"544681"
$arr[94,4] = [double]"2.5832514793933115e-08"
$arr[95,0] = [double]"23"
$arr[95,1] = "MD"
$arr[95,2] = [double]"251.79917907714844"
$arr[95,3] = [double]"544681"
$arr[95,4] = [double]"6.5919301128758434e-09"
$arr[96,0] = [double]"24"
$arr[96,1] = "MD"
$arr[96,2] = [double]"1"
$arr[96,3] = [double]"544681"
$arr[96,4] = [double]"2.5409620588168202e-11"
$arr[97,0] = [double]"25"
$arr[97,1] = "MD"
$arr[97,2] = [double]"1.5880948305130005"
$arr[97,3] = [double]"544681"
$arr[97,4] = [double]"8.8435474221437715e-11"
$arr[98,0] = [double]"26"
$arr[98,1] = "MD"
$arr[98,2] = [double]"575.76568603515625"
$arr[98,3] = [double]"544681"
$arr[98,4] = [double]"3.1763846664034645e-08"
$arr[99,0] = [double]"27"
$arr[99,1] = "MD"
$arr[99,2] = [double]"1"
$arr[99,3] = [double]"544681"
$arr[99,4] = [double]"1.201317517685041e-10"
$arr[100,0] = [double]"28"
$arr[100,1] = "MD"
$arr[100,2] = [double]"1"
$arr[100,3] = [double]"544681"
$arr[100,4] = [double]"1.1864852156318051e-10"
$arr[101,0] = [double]"29"
$arr[101,1] = "MD"
$arr[101,2] = [double]"1"
$arr[101,3] = [double]"544681"
$arr[101,4] = [double]"3.4642991253441835e-10"
$arr[102,0] = [double]"30"
$arr[102,1] = "MD"
$arr[102,2] = [double]"1"
$arr[102,3] = [double]"544681"
$arr[102,4] = [double]"8.8253360175372109e-10"
$arr[103,0] = [double]"31"
$arr[103,1] = "MD"
$arr[103,2] = [double]"1"
$arr[103,3] = [double]"544681"
$arr[103,4] = [double]"1.3238005136528841e-09"
$arr[104,0] = [double]"6"
$arr[104,1] = "NC"
$arr[104,2] = [double]"180.49602699279785"
$arr[104,3] = [double]"26646"
$arr[104,4] = [double]"1.4982975926614017e-07"
$arr[105,0] = [double]"7"
$arr[105,1] = "NC"
$arr[105,2] = [double]"152.10235595703125"
$arr[105,3] = [double]"26646"
$arr[105,4] = [double]"1.0286512086565835e-08"
$arr[106,0] = [double]"8"
$arr[106,1] = "NC"
$arr[106,2] = [double]"168.40696716308594"
$arr[106,3] = [double]"26646"
$arr[106,4] = [double]"9.4619867496703591e-09"
$arr[107,0] = [double]"9"
$arr[107,1] = "NC"
$arr[107,2] = [double]"218.52995300292969"
$arr[107,3] = [double]"26646"
$arr[107,4] = [double]"1.0841876729728028e-08"
$arr[108,0] = [double]"10"
$arr[108,1] = "NC"
$arr[108,2] = [double]"391.04666137695313"
$arr[108,3] = [double]"26646"
$arr[108,4] = [double]"4.9568402715749471e-08"
$arr[109,0] = [double]"11"
$arr[109,1] = "NC"
$arr[109,2] = [double]"679.3365478515625"
$arr[109,3] = [double]"26646"
$arr[109,4] = [double]"8.959761288451773e-08"
$arr[110,0] = [double]"12"
$arr[110,1] = "NC"
$arr[110,2] = [double]"428.40313720703125"
$arr[110,3] = [double]"26646"
$arr[110,4] = [double]"7.6333563470143417e-08"
$arr[111,0] = [double]"13"
$arr[111,1] = "NC"
$arr[111,2] = [double]"231.57841491699219"
$arr[111,3] = [double]"26646"
$arr[111,4] = [double]"3.0673472650732947e-08"
$arr[112,0] = [double]"14"
$arr[112,1] = "NC"
$arr[112,2] = [double]"95.509941101074219"
$arr[112,3] = [double]"26646"
$arr[112,4] = [double]"9.4166345832036313e-09"
$arr[113,0] = [double]"15"
$arr[113,1] = "NC"
$arr[113,2] = [double]"12171.7392578125"
$arr[113,3] = [double]"26646"
$arr[113,4] = [double]"6.9099161237318185e-07"
$arr[114,0] = [double]"16"
$arr[114,1] = "NC"
$arr[114,2] = [double]"6289.10595703125"
$arr[114,3] = [double]"26646"
$arr[114,4] = [double]"7.3007987566597876e-07"
$arr[115,0] = [double]"17"
$arr[115,1] = "NC"
$arr[115,2] = [double]"6330.93896484375"
$arr[115,3] = [double]"26646"
$arr[115,4] = [double]"8.0649090250517474e-07"
$arr[116,0] = [double]"18"
$arr[116,1] = "NC"
$arr[116,2] = [double]"3262.856201171875"
$arr[116,3] = [double]"26646"
$arr[116,4] = [double]"7.287360404006904e-07"
$arr[117,0] = [double]"19"
$arr[117,1] = "NC"
$arr[117,2] = [double]"3283.679443359375"
$arr[117,3] = [double]"26646"
$arr[117,4] = [double]"7.6596120379690547e-07"
$arr[118,0] = [double]"20"
$arr[118,1] = "NC"
$arr[118,2] = [double]"1694.51806640625"
$arr[118,3] = [double]"26646"
$arr[118,4] = [double]"4.8234448968287325e-07"
$arr[119,0] = [double]"21"
$arr[119,1] = "NC"
$arr[119,2] = [double]"366.25119018554688"
$arr[119,3] = [double]"26646"
$arr[119,4] = [double]"1.947375523059236e-07"
$arr[120,0] = [double]"22"
$arr[120,1] = "NC"
$arr[120,2] = [double]"262.08477783203125"
$arr[120,3] = [double]"26646"
$arr[120,4] = [double]"9.9614887005827768e-08"
$arr[121,0] = [double]"23"
$arr[121,1] = "NC"
$arr[121,2] = [double]"47.500885009765625"
$arr[121,3] = [double]"26646"
$arr[121,4] = [double]"2.5419685911742818e-08"
$arr[122,0] = [double]"24"
$arr[122,1] = "NC"
$arr[122,2] = [double]"1"
$arr[122,3] = [double]"26646"
$arr[122,4] = [double]"5.1940768353020417e-10"
$arr[123,0] = [double]"25"
$arr[123,1] = "NC"
$arr[123,2] = [double]"0.29958760738372803"
$arr[123,3] = [double]"26646"
$arr[123,4] = [double]"3.4102334844909876e-10"
$arr[124,0] = [double]"26"
$arr[124,1] = "NC"
$arr[124,2] = [double]"108.61583709716797"
$arr[124,3] = [double]"26646"
$arr[124,4] = [double]"1.2248719372109917e-07"
$arr[125,0] = [double]"27"
$arr[125,1] = "NC"
$arr[125,2] = [double]"1"
$arr[125,3] = [double]"26646"
$arr[125,4] = [double]"2.4556585653812135e-09"
$arr[126,0] = [double]"28"
$arr[126,1] = "NC"
$arr[126,2] = [double]"1"
$arr[126,3] = [double]"26646"
$arr[126,4] = [double]"2.4253394848017251e-09"
$arr[127,0] = [double]"29"
$arr[127,1] = "NC"
$arr[127,2] = [double]"1"
$arr[127,3] = [double]"26646"
$arr[127,4] = [double]"7.0815056041340085e-09"
$arr[128,0] = [double]"30"
$arr[128,1] = "NC"
$arr[128,2] = [double]"1"
$arr[128,3] = [double]"26646"
$arr[128,4] = [double]"1.8040203997315984e-08"
$arr[129,0] = [double]"31"
$arr[129,1] = "NC"
$arr[129,2] = [double]"1"
$arr[129,3] = [double]"26646"
$arr[129,4] = [double]"2.7060307772330816e-08"
$arr[130,0] = [double]"6"
$arr[130,1] = "NJ"
$arr[130,2] = [double]"13443.02294921875"
$arr[130,3] = [double]"4579196"
$arr[130,4] = [double]"1.2946517635725741e-08"
$arr[131,0] = [double]"7"
$arr[131,1] = "NJ"
$arr[131,2] = [double]"9091.7099609375"
$arr[131,3] = [double]"4579196"
$arr[131,4] = [double]"7.1335115592319198e-10"
$arr[132,0] = [double]"8"
$arr[132,1] = "NJ"
$arr[132,2] = [double]"18183.419921875"
$arr[132,3] = [double]"4579196"
$arr[132,4] = [double]"1.1852867576322978e-09"
$arr[133,0] = [double]"9"
$arr[133,1] = "NJ"
$arr[133,2] = [double]"68187.8203125"
$arr[133,3] = [double]"4579196"
$arr[133,4] = [double]"3.9248755356879883e-09"
$arr[134,0] = [double]"10"
$arr[134,1] = "NJ"
$arr[134,2] = [double]"388670.59375"
$arr[134,3] = [double]"4579196"
$arr[134,4] = [double]"5.715886786106239e-08"
$arr[135,0] = [double]"11"
$arr[135,1] = "NJ"
$arr[135,2] = [double]"892014.625"
$arr[135,3] = [double]"4579196"
$arr[135,4] = [double]"1.3649257368797407e-07"
$arr[136,0] = [double]"12"
$arr[136,1] = "NJ"
$arr[136,2] = [double]"1672874.625"
$arr[136,3] = [double]"4579196"
$arr[136,4] = [double]"3.4582140528982563e-07"
$arr[137,0] = [double]"13"
$arr[137,1] = "NJ"
$arr[137,2] = [double]"1722879.25"
$arr[137,3] = [double]"4579196"
$arr[137,4] = [double]"2.6475569825379353e-07"
$arr[138,0] = [double]"14"
$arr[138,1] = "NJ"
$arr[138,2] = [double]"2179945.75"
$arr[138,3] = [double]"4579196"
$arr[138,4] = [double]"2.493551960469631e-07"
$arr[139,0] = [double]"15"
$arr[139,1] = "NJ"
$arr[139,2] = [double]"2146865"
$arr[139,3] = [double]"4579196"
$arr[139,4] = [double]"1.4140033499643323e-07"
$arr[140,0] = [double]"16"
$arr[140,1] = "NJ"
$arr[140,2] = [double]"2126356.5"
$arr[140,3] = [double]"4579196"
$arr[140,4] = [double]"2.8638029903049755e-07"
$arr[141,0] = [double]"17"
$arr[141,1] = "NJ"
$arr[141,2] = [double]"1803709.375"
$arr[141,3] = [double]"4579196"
$arr[141,4] = [double]"2.6657747298486356e-07"
$arr[142,0] = [double]"18"
$arr[142,1] = "NJ"
$arr[142,2] = [double]"557451.125"
$arr[142,3] = [double]"4579196"
$arr[142,4] = [double]"1.4444572116190102e-07"
$arr[143,0] = [double]"19"
$arr[143,1] = "NJ"
$arr[143,2] = [double]"222423.40625"
$arr[143,3] = [double]"4579196"
$arr[143,4] = [double]"6.0193848128164973e-08"
$arr[144,0] = [double]"20"
$arr[144,1] = "NJ"
$arr[144,2] = [double]"130643.3671875"
$arr[144,3] = [double]"4579196"
$arr[144,4] = [double]"4.3144360972746654e-08"
$arr[145,0] = [double]"21"
$arr[145,1] = "NJ"
$arr[145,2] = [double]"78947.546875"
$arr[145,3] = [double]"4579196"
$arr[145,4] = [double]"4.8700663057843485e-08"
$arr[146,0] = [double]"22"
$arr[146,1] = "NJ"
$arr[146,2] = [double]"53178.76953125"
$arr[146,3] = [double]"4579196"
$arr[146,4] = [double]"2.3450185793194578e-08"
$arr[147,0] = [double]"23"
$arr[147,1] = "NJ"
$arr[147,2] = [double]"18301.86328125"
$arr[147,3] = [double]"4579196"
$arr[147,4] = [double]"1.1362905283363034e-08"
$arr[148,0] = [double]"24"
$arr[148,1] = "NJ"
$arr[148,2] = [double]"33230.90625"
$arr[148,3] = [double]"4579196"
$arr[148,4] = [double]"2.0025167302151203e-08"
$arr[149,0] = [double]"25"
$arr[149,1] = "NJ"
$arr[149,2] = [double]"28838.1875"
$arr[149,3] = [double]"4579196"
$arr[149,4] = [double]"3.8084987608044685e-08"
$arr[150,0] = [double]"26"
$arr[150,1] = "NJ"
$arr[150,2] = [double]"5879.18310546875"
$arr[150,3] = [double]"4579196"
$arr[150,4] = [double]"7.6920141367509132e-09"
$arr[151,0] = [double]"27"
$arr[151,1] = "NJ"
$arr[151,2] = [double]"314.47021484375"
$arr[151,3] = [double]"4579196"
$arr[151,4] = [double]"8.9592794294546252e-10"
$arr[152,0] = [double]"28"
$arr[152,1] = "NJ"
$arr[152,2] = [double]"1"
$arr[152,3] = [double]"4579196"
$arr[152,4] = [double]"2.8138316329751056e-12"
$arr[153,0] = [double]"29"
$arr[153,1] = "NJ"
$arr[153,2] = [double]"1572.2359619140625"
$arr[153,3] = [double]"4579196"
$arr[153,4] = [double]"1.2917215741481414e-08"
$arr[154,0] = [double]"30"
$arr[154,1] = "NJ"
$arr[154,2] = [double]"1662.5887451171875"
$arr[154,3] = [double]"4579196"
$arr[154,4] = [double]"3.4797807302311412e-08"
$arr[155,0] = [double]"31"
$arr[155,1] = "NJ"
$arr[155,2] = [double]"1662.5887451171875"
$arr[155,3] = [double]"4579196"
$arr[155,4] = [double]"5.2196710953467118e-08"
$arr[156,0] = [double]"6"
$arr[156,1] = "NY"
$arr[156,2] = [double]"4344.2998046875"
$arr[156,3] = [double]"3474730"
$arr[156,4] = [double]"5.9955231890285177e-09"
$arr[157,0] = [double]"7"
$arr[157,1] = "NY"
$arr[157,2] = [double]"20007.62109375"
$arr[157,3] = [double]"3474730"
$arr[157,4] = [double]"2.2495960649848712e-09"
$arr[158,0] = [double]"8"
$arr[158,1] = "NY"
$arr[158,2] = [double]"28899.896484375"
$arr[158,3] = [double]"3474730"
$arr[158,4] = [double]"2.6995752300251752e-09"
$arr[159,0] = [double]"9"
$arr[159,1] = "NY"
$arr[159,2] = [double]"46684.4453125"
$arr[159,3] = [double]"3474730"
$arr[159,4] = [double]"3.8507255162301135e-09"
$arr[160,0] = [double]"10"
$arr[160,1] = "NY"
$arr[160,2] = [double]"157837.875"
$arr[160,3] = [double]"3474730"
$arr[160,4] = [double]"3.326323039232193e-08"
$arr[161,0] = [double]"11"
$arr[161,1] = "NY"
$arr[161,2] = [double]"549351.25"
$arr[161,3] = [double]"3474730"
$arr[161,4] = [double]"1.2045877895161539e-07"
$arr[162,0] = [double]"12"
$arr[162,1] = "NY"
$arr[162,2] = [double]"936842.6875"
$arr[162,3] = [double]"3474730"
$arr[162,4] = [double]"2.7752776077250019e-07"
$arr[163,0] = [double]"13"
$arr[163,1] = "NY"
$arr[163,2] = [double]"1302913.375"
$arr[163,3] = [double]"3474730"
$arr[163,4] = [double]"2.869176682906982e-07"
$arr[164,0] = [double]"14"
$arr[164,1] = "NY"
$arr[164,2] = [double]"1570005.25"
$arr[164,3] = [double]"3474730"
$arr[164,4] = [double]"2.5735056397024891e-07"
$arr[165,0] = [double]"15"
$arr[165,1] = "NY"
$arr[165,2] = [double]"1674901.5"
$arr[165,3] = [double]"3474730"
$arr[165,4] = [double]"1.5808340947387478e-07"
$arr[166,0] = [double]"16"
$arr[166,1] = "NY"
$arr[166,2] = [double]"1351896.875"
$arr[166,3] = [double]"3474730"
$arr[166,4] = [double]"2.6091672111761e-07"
$arr[167,0] = [double]"17"
$arr[167,1] = "NY"
$arr[167,2] = [double]"856904.5625"
$arr[167,3] = [double]"3474730"
$arr[167,4] = [double]"1.814849639458771e-07"
$arr[168,0] = [double]"18"
$arr[168,1] = "NY"
$arr[168,2] = [double]"469480"
$arr[168,3] = [double]"3474730"
$arr[168,4] = [double]"1.7432765275771089e-07"
$arr[169,0] = [double]"19"
$arr[169,1] = "NY"
$arr[169,2] = [double]"165101.734375"
$arr[169,3] = [double]"3474730"
$arr[169,4] = [double]"6.4028675694771664e-08"
$arr[170,0] = [double]"20"
$arr[170,1] = "NY"
$arr[170,2] = [double]"96088.796875"
$arr[170,3] = [double]"3474730"
$arr[170,4] = [double]"4.5473743881530027e-08"
$arr[171,0] = [double]"21"
$arr[171,1] = "NY"
$arr[171,2] = [double]"100041.7109375"
$arr[171,3] = [double]"3474730"
$arr[171,4] = [double]"8.8435903933259397e-08"
$arr[172,0] = [double]"22"
$arr[172,1] = "NY"
$arr[172,2] = [double]"75425.0625"
$arr[172,3] = [double]"3474730"
$arr[172,4] = [double]"4.7662297220085748e-08"
$arr[173,0] = [double]"23"
$arr[173,1] = "NY"
$arr[173,2] = [double]"54351.703125"
$arr[173,3] = [double]"3474730"
$arr[173,4] = [double]"4.8356906034996427e-08"
$arr[174,0] = [double]"24"
$arr[174,1] = "NY"
$arr[174,2] = [double]"30674.642578125"
$arr[174,3] = [double]"3474730"
$arr[174,4] = [double]"2.6488949700365083e-08"
$arr[175,0] = [double]"25"
$arr[175,1] = "NY"
$arr[175,2] = [double]"28858.61328125"
$arr[175,3] = [double]"3474730"
$arr[175,4] = [double]"5.4615082234477086e-08"
$arr[176,0] = [double]"26"
$arr[176,1] = "NY"
$arr[176,2] = [double]"20119.1953125"
$arr[176,3] = [double]"3474730"
$arr[176,4] = [double]"3.7721154200198725e-08"
$arr[177,0] = [double]"27"
$arr[177,1] = "NY"
$arr[177,2] = [double]"10125.1298828125"
$arr[177,3] = [double]"3474730"
$arr[177,4] = [double]"4.1337617773251623e-08"
$arr[178,0] = [double]"28"
$arr[178,1] = "NY"
$arr[178,2] = [double]"5227.9228515625"
$arr[178,3] = [double]"3474730"
$arr[178,4] = [double]"2.1080387213601171e-08"
$arr[179,0] = [double]"29"
$arr[179,1] = "NY"
$arr[179,2] = [double]"1749.8079833984375"
$arr[179,3] = [double]"3474730"
$arr[179,4] = [double]"2.0601216732529792e-08"
$arr[180,0] = [double]"30"
$arr[180,1] = "NY"
$arr[180,2] = [double]"3758.697998046875"
$arr[180,3] = [double]"3474730"
$arr[180,4] = [double]"1.1273421307578246e-07"
$arr[181,0] = [double]"31"
$arr[181,1] = "NY"
$arr[181,2] = [double]"1381.36962890625"
$arr[181,3] = [double]"3474730"
$arr[181,4] = [double]"6.2146902735094045e-08"
$arr[182,0] = [double]"6"
$arr[182,1] = "RI"
$arr[182,2] = [double]"1"
$arr[182,3] = [double]"472523"
$arr[182,4] = [double]"1.0148583362568075e-11"
$arr[183,0] = [double]"7"
$arr[183,1] = "RI"
$arr[183,2] = [double]"1"
$arr[183,3] = [double]"472523"
$arr[183,4] = [double]"8.2681279357788018e-13"
$arr[184,0] = [double]"8"
$arr[184,1] = "RI"
$arr[184,2] = [double]"3634.8701171875"
$arr[184,3] = [double]"472523"
$arr[184,4] = [double]"2.4968138667702533e-09"
$arr[185,0] = [double]"9"
$arr[185,1] = "RI"
$arr[185,2] = [double]"10904.6103515625"
$arr[185,3] = [double]"472523"
$arr[185,4] = [double]"6.6142198384966377e-09"
$arr[186,0] = [double]"10"
$arr[186,1] = "RI"
$arr[186,2] = [double]"9692.9873046875"
$arr[186,3] = [double]"472523"
$arr[186,4] = [double]"1.5021349497601477e-08"
$arr[187,0] = [double]"11"
$arr[187,1] = "RI"
$arr[187,2] = [double]"33925.453125"
$arr[187,3] = [double]"472523"
$arr[187,4] = [double]"5.4703122032151441e-08"
$arr[188,0] = [double]"12"
$arr[188,1] = "RI"
$arr[188,2] = [double]"100393"
$arr[188,3] = [double]"472523"
$arr[188,4] = [double]"2.186962291261807e-07"
$arr[189,0] = [double]"13"
$arr[189,1] = "RI"
$arr[189,2] = [double]"134507.125"
$arr[189,3] = [double]"472523"
$arr[189,4] = [double]"2.1781367820494779e-07"
$arr[190,0] = [double]"14"
$arr[190,1] = "RI"
$arr[190,2] = [double]"144514.234375"
$arr[190,3] = [double]"472523"
$arr[190,4] = [double]"1.7419381492800312e-07"
$arr[191,0] = [double]"15"
$arr[191,1] = "RI"
$arr[191,2] = [double]"233576.234375"
$arr[191,3] = [double]"472523"
$arr[191,4] = [double]"1.6211521369768889e-07"
$arr[192,0] = [double]"16"
$arr[192,1] = "RI"
$arr[192,2] = [double]"243603.46875"
$arr[192,3] = [double]"472523"
$arr[192,4] = [double]"3.4573227480905189e-07"
$arr[193,0] = [double]"17"
$arr[193,1] = "RI"
$arr[193,2] = [double]"254660.453125"
$arr[193,3] = [double]"472523"
$arr[193,4] = [double]"3.9661378536948177e-07"
$arr[194,0] = [double]"18"
$arr[194,1] = "RI"
$arr[194,2] = [double]"107529.3984375"
$arr[194,3] = [double]"472523"
$arr[194,4] = [double]"2.9361248721215816e-07"
$arr[195,0] = [double]"19"
$arr[195,1] = "RI"
$arr[195,2] = [double]"58113.5859375"
$arr[195,3] = [double]"472523"
$arr[195,4] = [double]"1.6572906247347419e-07"
$arr[196,0] = [double]"20"
$arr[196,1] = "RI"
$arr[196,2] = [double]"67729.8125"
$arr[196,3] = [double]"472523"
$arr[196,4] = [double]"2.3570345319967601e-07"
$arr[197,0] = [double]"21"
$arr[197,1] = "RI"
$arr[197,2] = [double]"58950.8203125"
$arr[197,3] = [double]"472523"
$arr[197,4] = [double]"3.8320879980346945e-07"
$arr[198,0] = [double]"22"
$arr[198,1] = "RI"
$arr[198,2] = [double]"17435"
$arr[198,3] = [double]"472523"
$arr[198,4] = [double]"8.1017596187393792e-08"
$arr[199,0] = [double]"23"
$arr[199,1] = "RI"
$arr[199,2] = [double]"13030.740234375"
$arr[199,3] = [double]"472523"
$arr[199,4] = [double]"8.5253567760901205e-08"
$arr[200,0] = [double]"24"
$arr[200,1] = "RI"
$arr[200,2] = [double]"17672.072265625"
$arr[200,3] = [double]"472523"
$arr[200,4] = [double]"1.1222014251188739e-07"
$arr[201,0] = [double]"25"
$arr[201,1] = "RI"
$arr[201,2] = [double]"3065.32470703125"
$arr[201,3] = [double]"472523"
$arr[201,4] = [double]"4.2659102916786651e-08"
$arr[202,0] = [double]"26"
$arr[202,1] = "RI"
$arr[202,2] = [double]"4138.96240234375"
$arr[202,3] = [double]"472523"
$arr[202,4] = [double]"5.7064230674086502e-08"
$arr[203,0] = [double]"27"
$arr[203,1] = "RI"
$arr[203,2] = [double]"11138.72265625"
$arr[203,3] = [double]"472523"
$arr[203,4] = [double]"3.3440929314565437e-07"
$arr[204,0] = [double]"28"
$arr[204,1] = "RI"
$arr[204,2] = [double]"575.0059814453125"
$arr[204,3] = [double]"472523"
$arr[204,4] = [double]"1.7049822886860966e-08"
$arr[205,0] = [double]"29"
$arr[205,1] = "RI"
$arr[205,2] = [double]"4134.5166015625"
$arr[205,3] = [double]"472523"
$arr[205,4] = [double]"3.5795244457403896e-07"
$arr[206,0] = [double]"30"
$arr[206,1] = "RI"
$arr[206,2] = [double]"16.977466583251953"
$arr[206,3] = [double]"472523"
$arr[206,4] = [double]"3.7444607414727216e-09"
$arr[207,0] = [double]"31"
$arr[207,1] = "RI"
$arr[207,2] = [double]"1"
$arr[207,3] = [double]"472523"
$arr[207,4] = [double]"3.3083213946127898e-10"
$arr[208,0] = [double]"6"
$arr[208,1] = "VA"
$arr[208,2] = [double]"42157.80029296875"
$arr[208,3] = [double]"785624"
$arr[208,4] = [double]"1.1869315130752511e-06"
$arr[209,0] = [double]"7"
$arr[209,1] = "VA"
$arr[209,2] = [double]"72246.0390625"
$arr[209,3] = [double]"785624"
$arr[209,4] = [double]"1.6571564742662304e-07"
$arr[210,0] = [double]"8"
$arr[210,1] = "VA"
$arr[210,2] = [double]"79990.453125"
$arr[210,3] = [double]"785624"
$arr[210,4] = [double]"1.52432548361503e-07"
$arr[211,0] = [double]"9"
$arr[211,1] = "VA"
$arr[211,2] = [double]"103798.0234375"
$arr[211,3] = [double]"785624"
$arr[211,4] = [double]"1.7466257418163877e-07"
$arr[212,0] = [double]"10"
$arr[212,1] = "VA"
$arr[212,2] = [double]"185740.53125"
$arr[212,3] = [double]"785624"
$arr[212,4] = [double]"7.9854669365886366e-07"
$arr[213,0] = [double]"11"
$arr[213,1] = "VA"
$arr[213,2] = [double]"322673.34375"
$arr[213,3] = [double]"785624"
$arr[213,4] = [double]"1.4434170907406951e-06"
$arr[214,0] = [double]"12"
$arr[214,1] = "VA"
$arr[214,2] = [double]"203484.234375"
$arr[214,3] = [double]"785624"
$arr[214,4] = [double]"1.2297333569222246e-06"
$arr[215,0] = [double]"13"
$arr[215,1] = "VA"
$arr[215,2] = [double]"109995.8203125"
$arr[215,3] = [double]"785624"
$arr[215,4] = [double]"4.9414950353821041e-07"
$arr[216,0] = [double]"14"
$arr[216,1] = "VA"
$arr[216,2] = [double]"89069.2734375"
$arr[216,3] = [double]"785624"
$arr[216,4] = [double]"2.9784638400087715e-07"
$arr[217,0] = [double]"15"
$arr[217,1] = "VA"
$arr[217,2] = [double]"94129.8828125"
$arr[217,3] = [double]"785624"
$arr[217,4] = [double]"1.8124454470580531e-07"
$arr[218,0] = [double]"16"
$arr[218,1] = "VA"
$arr[218,2] = [double]"48507.75"
$arr[218,3] = [double]"785624"
$arr[218,4] = [double]"1.9098968095931923e-07"
$arr[219,0] = [double]"17"
$arr[219,1] = "VA"
$arr[219,2] = [double]"68175.4453125"
$arr[219,3] = [double]"785624"
$arr[219,4] = [double]"2.945619144156808e-07"
$arr[220,0] = [double]"18"
$arr[220,1] = "VA"
$arr[220,2] = [double]"35136.44140625"
$arr[220,3] = [double]"785624"
$arr[220,4] = [double]"2.6616280024427397e-07"
$arr[221,0] = [double]"19"
$arr[221,1] = "VA"
$arr[221,2] = [double]"35360.6796875"
$arr[221,3] = [double]"785624"
$arr[221,4] = [double]"2.7975889338449633e-07"
$arr[222,0] = [double]"20"
$arr[222,1] = "VA"
$arr[222,2] = [double]"18247.61328125"
$arr[222,3] = [double]"785624"
$arr[222,4] = [double]"1.7617101377709332e-07"
$arr[223,0] = [double]"21"
$arr[223,1] = "VA"
$arr[223,2] = [double]"3944.017822265625"
$arr[223,3] = [double]"785624"
$arr[223,4] = [double]"7.1125739964372769e-08"
$arr[224,0] = [double]"22"
$arr[224,1] = "VA"
$arr[224,2] = [double]"2822.2900390625"
$arr[224,3] = [double]"785624"
$arr[224,4] = [double]"3.6383237755899245e-08"
$arr[225,0] = [double]"23"
$arr[225,1] = "VA"
$arr[225,2] = [double]"511.51873779296875"
$arr[225,3] = [double]"785624"
$arr[225,4] = [double]"9.2842604715315247e-09"
$arr[226,0] = [double]"24"
$arr[226,1] = "VA"
$arr[226,2] = [double]"1"
$arr[226,3] = [double]"785624"
$arr[226,4] = [double]"1.761674486844278e-11"
$arr[227,0] = [double]"25"
$arr[227,1] = "VA"
$arr[227,2] = [double]"3.2261433601379395"
$arr[227,3] = [double]"785624"
$arr[227,4] = [double]"1.245550190764888e-10"
$arr[228,0] = [double]"26"
$arr[228,1] = "VA"
$arr[228,2] = [double]"1169.64208984375"
$arr[228,3] = [double]"785624"
$arr[228,4] = [double]"4.47370958056581e-08"
$arr[229,0] = [double]"27"
$arr[229,1] = "VA"
$arr[229,2] = [double]"1"
$arr[229,3] = [double]"785624"
$arr[229,4] = [double]"8.3288542729320625e-11"
$arr[230,0] = [double]"28"
$arr[230,1] = "VA"
$arr[230,2] = [double]"1"
$arr[230,3] = [double]"785624"
$arr[230,4] = [double]"8.2260212530549381e-11"
$arr[231,0] = [double]"29"
$arr[231,1] = "VA"
$arr[231,2] = [double]"1"
$arr[231,3] = [double]"785624"
$arr[231,4] = [double]"2.4018334543463027e-10"
$arr[232,0] = [double]"30"
$arr[232,1] = "VA"
$arr[232,2] = [double]"1"
$arr[232,3] = [double]"785624"
$arr[232,4] = [double]"6.1186938848933892e-10"
$arr[233,0] = [double]"31"
$arr[233,1] = "VA"
$arr[233,2] = [double]"1"
$arr[233,3] = [double]"785624"
$arr[233,4] = [double]"9.1780416600073522e-10"

$ws.Range("A2:E235").Value = $arr

$fmtA = $ws.Cells.Item(2,1).NumberFormat
$fmtC = $ws.Cells.Item(2,3).NumberFormat
$fmtD = $ws.Cells.Item(2,4).NumberFormat
$fmtE = $ws.Cells.Item(2,5).NumberFormat
$ws.Range("A231:A235").NumberFormat = $fmtA
$ws.Range("C231:C235").NumberFormat = $fmtC
$ws.Range("D231:D235").NumberFormat = $fmtD
$ws.Range("E231:E235").NumberFormat = $fmtE
